$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new "TestSuite" sheet right after "addCustomerTest"
#    (so the tab order becomes addCustomerTest, TestSuite, openAccountTest)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("addCustomerTest")
$suite = $wb.Worksheets.Add($null, $ws1)
$suite.Name = "TestSuite"

# Re-resolve sheets by name *after* the insert: adding a sheet shifts the
# underlying positions, and handles grabbed beforehand can go stale.
$ws1 = $wb.Worksheets.Item("addCustomerTest")
$ws3 = $wb.Worksheets.Item("openAccountTest")
$suite = $wb.Worksheets.Item("TestSuite")

# ---------------------------------------------------------------------
# 2. addCustomerTest: append 3 more customer rows (rows 3-5)
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "Radha"
$ws1.Range("B3").Value = "Madhusudan"
$ws1.Range("C3").Value = 354357
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Radha"
$ws1.Range("B4").Value = "Rasbihari"
$ws1.Range("C4").Value = 354358
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "Govinda"
$ws1.Range("B5").Value = "Gopal"
$ws1.Range("C5").Value = 354359
$ws1.Range("D5").Value = "Customer added successfully"

$ws1.Range("A1:D5").WrapText = $true
$ws1.Rows.Item(1).RowHeight = 30
$ws1.Rows.Item(2).RowHeight = 30
$ws1.Rows.Item(3).RowHeight = 30
$ws1.Rows.Item(4).RowHeight = 30
$ws1.Rows.Item(5).RowHeight = 30
$ws1.Columns.Item(4).ColumnWidth = 26.28515625

[void]$ws1.Range("F4").Select()

# ---------------------------------------------------------------------
# 3. TestSuite: new TCID / RUNMODE control table
# ---------------------------------------------------------------------
$suite.Range("A1").Value = "TCID"
$suite.Range("B1").Value = "RUNMODE"

$suite.Range("A2").Value = "LoginTest"
$suite.Range("B2").Value = "Y"

$suite.Range("A3").Value = "AddCustomerTest"
$suite.Range("B3").Value = "Y"

$suite.Range("A4").Value = "OpenAccountTest"
$suite.Range("B4").Value = "N"

$suite.Range("A1:B4").WrapText = $true
$suite.Rows.Item(2).RowHeight = 30
$suite.Rows.Item(3).RowHeight = 45
$suite.Rows.Item(4).RowHeight = 30
$suite.Columns.Item(2).ColumnWidth = 10.42578125

# ---------------------------------------------------------------------
# 4. openAccountTest: replace old content with the Customer/Currency/
#    AlertText table used to validate the "open account" flow.
# ---------------------------------------------------------------------
$ws3.Cells.Clear()

$ws3.Range("A1").Value = "Customer"
$ws3.Range("B1").Value = "Currency"
$ws3.Range("C1").Value = "AlertText"

$ws3.Range("A2").Value = "Radha Raman"
$ws3.Range("B2").Value = "Rupee"
$ws3.Range("C2").Value = "Account created successfully with account Number "

$ws3.Range("A1:C2").WrapText = $true
$ws3.Rows.Item(1).RowHeight = 30
$ws3.Rows.Item(2).RowHeight = 90
$ws3.Columns.Item(1).ColumnWidth = 13.85546875

[void]$ws3.Range("F7").Select()

# ---------------------------------------------------------------------
# 5. Make TestSuite the active tab / sheet, with B5 selected
# ---------------------------------------------------------------------
[void]$suite.Range("B5").Select()
$suite.Activate()

Write-Host "edit complete"
